$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 12 ("Legislature") - this shifts all subsequent rows up by one,
# restoring the "rural communities" interactions that belong to the rows below
# back onto their correct labels.
$ws.Rows.Item(12).Delete()
